# -------------------------------------------------------------------------
# Applies the "DistOriginal dentro dentro da Grafo" commit:
#   1. "30/10 às 9:10" -> "23/10 às 09h10min", appends a new sentence
#      listing classes added (several italicized class names), and moves
#      the _GoBack bookmark into this paragraph (between "Git" and "Hub").
#   2. "05/11: criação das classes ..." -> "05/11: adaptação das classes
#      ... e criação da classe Cidade; ..." (Cidade + RelativeLayout now
#      italic).
#   3. "05/11: por alguma razão ... RelativeLayout ..." merges the trailing
#      "Provavelmente..." sentence into the paragraph and italicizes
#      RelativeLayout.
#   4. "Resolução: ... reconhecido." -> "... reconhecidos." and the
#      old _GoBack bookmark (that used to sit here) is removed.
# -------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0. Drop the pre-existing _GoBack bookmark (it will be re-created later
#    in its new location, paragraph 12).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1. Paragraph 12: "30/10 às 9:10: início do desenvolvimento ..."
# ---------------------------------------------------------------------
$old1 = "30/10 às 9:10: início do desenvolvimento do projeto. Criamos o repositório no GitHub e optamos pelo Xamarin – programação Android no Visual Studio, na linguagem C#."
$new1 = "23/10 às 09h10min: início do desenvolvimento do projeto. Criamos o repositório no GitHub e optamos pelo Xamarin – programação Android no Visual Studio, na linguagem C#. Adição das classes BucketHash, No, IStack, Pilha, PilhaVaziaException, Grafo e Vertice, que desenvolvemos anteriormente na matéria."
$d.Paragraphs(12).Range.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

$italicWords1 = @("BucketHash", "No", "IStack", "Pilha", "PilhaVaziaException", "Grafo", "Vertice")
$searchStart = $d.Paragraphs(12).Range.Start
$paraEnd = $d.Paragraphs(12).Range.End
foreach ($w in $italicWords1) {
    $rng = $d.Range($searchStart, $paraEnd)
    $found = $rng.Find.Execute($w, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Italic = 1
        $searchStart = $rng.End
        $paraEnd = $d.Paragraphs(12).Range.End
    }
}

# Re-insert the _GoBack bookmark between "Git" and "Hub".
$bmFindRng = $d.Paragraphs(12).Range
$bmFindRng.Find.Execute("repositório no Git", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPos = $bmFindRng.End
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 2. Paragraph 13: "05/11: criação das classes necessárias ..."
# ---------------------------------------------------------------------
$old2 = "05/11: criação das classes necessárias para o projeto, com base nas classes desenvolvidas previamente na matéria; começo do design da página inicial usando o RelativeLayout."
$new2 = "05/11: adaptação das classes necessárias para o projeto e criação da classe Cidade; começo do design da página inicial usando o RelativeLayout."
$d.Paragraphs(13).Range.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

$italicWords2 = @("Cidade", "RelativeLayout")
$searchStart = $d.Paragraphs(13).Range.Start
$paraEnd = $d.Paragraphs(13).Range.End
foreach ($w in $italicWords2) {
    $rng = $d.Range($searchStart, $paraEnd)
    $found = $rng.Find.Execute($w, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Italic = 1
        $searchStart = $rng.End
        $paraEnd = $d.Paragraphs(13).Range.End
    }
}

# ---------------------------------------------------------------------
# 3. Paragraph 16: "05/11: por alguma razão ... RelativeLayout ..."
#    (text is unchanged, but re-typing it collapses the stray proofErr
#    marks and lets us italicize RelativeLayout cleanly)
# ---------------------------------------------------------------------
$old3 = "05/11: por alguma razão, os atributos que seriam usados para posicionar os elementos no RelativeLayout não estão sendo reconhecidos. Provavelmente está faltando baixar algo no Visual Studio."
$d.Paragraphs(16).Range.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $old3, 2)

$rng3 = $d.Paragraphs(16).Range
$found3 = $rng3.Find.Execute("RelativeLayout", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Italic = 1
}

# ---------------------------------------------------------------------
# 4. Paragraph 17: "Resolução: ... reconhecido." -> "... reconhecidos."
# ---------------------------------------------------------------------
$old4 = "Resolução: atualizei o projeto, então, os atributos foram reconhecido."
$new4 = "Resolução: atualizei o projeto, então, os atributos foram reconhecidos."
$d.Paragraphs(17).Range.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
